# Update the FHIR StructureDefinition metadata/elements workbook from the
# "Alvearie"/IBM publication to the "LinuxForHealth" publication.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: Property / Value table ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/snapshot-age-in-years"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: element definition table ---
$elements = $wb.Worksheets.Item("Elements")

# Extension.url's Fixed Value mirrors the StructureDefinition URL above.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/snapshot-age-in-years"

# The top-level "Extension" row's Constraint(s) cell was a stray duplicate of
# the Extension.extension row's constraints; clear it so only
# Extension.extension (AI4) carries that text.
$elements.Range("AI2").Value = ""
